$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B (min), C (median), D (max) for rows 2-19 (PCB species table)
$data = @(
    @(0.77,1.05,1.15),
    @(0.37,0.63,0.93),
    @(0.12,0.3,0.72),
    @(0.54,0.8100000000000001,1.4),
    @(0.54,0.98,1.32),
    @(0.63,1.21,1.63),
    @(0.73,1.02,1.41),
    @(0.66,0.99,1.34),
    @(0.65,0.89,1.19),
    @(0.87,1.16,1.52),
    @(0.72,0.95,1.16),
    @(0.5,0.68,0.89),
    @(0.79,1.03,1.35),
    @(0.65,0.77,0.98),
    @(0.83,1.1,1.28),
    @(0.72,1.04,1.32),
    @(0.72,1.11,1.64),
    @(0.75,1.16,1.71)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 2).Value = $data[$i][0]
    $ws.Cells.Item($r, 3).Value = $data[$i][1]
    $ws.Cells.Item($r, 4).Value = $data[$i][2]
}
